$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Inscritos (E3) 18 -> 19
$ws.Range("E3").Value = 19

# Row 18: Inscritos (E18) 89 -> 90
$ws.Range("E18").Value = 90

# Row 31: Inscritos (E31) 2 -> 1, Pagos (F31) 1 -> 0, Inscricoes homologadas (H31) 1 -> 0
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 0
$ws.Range("H31").Value = 0

# Row 38: Inscritos (E38) 54 -> 55
$ws.Range("E38").Value = 55

# Row 63: Inscritos (E63) 22 -> 24
$ws.Range("E63").Value = 24

# Row 72: Inscritos (E72) 30 -> 31
$ws.Range("E72").Value = 31

# Row 73: Inscritos (E73) 22 -> 24
$ws.Range("E73").Value = 24

# Row 76: Inscritos (E76) 37 -> 38
$ws.Range("E76").Value = 38

# Row 86: Inscritos (E86) 1 -> 2, Pagos (F86) 0 -> 1, Inscricoes homologadas (H86) 0 -> 1
$ws.Range("E86").Value = 2
$ws.Range("F86").Value = 1
$ws.Range("H86").Value = 1

$wb.Save()
